$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3..30 down to 4..31.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly record.
$ws.Cells.Item(3, 1).Value = 4
$ws.Cells.Item(3, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(3, 3).Value = "Los Lagos"
$ws.Cells.Item(3, 4).Value = 45163
$ws.Cells.Item(3, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3, 5).Value = 10
$ws.Cells.Item(3, 6).Value = 100112035
$ws.Cells.Item(3, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 140
$ws.Cells.Item(3, 11).Value = 25000
$ws.Cells.Item(3, 12).Value = 25000
$ws.Cells.Item(3, 13).Value = 25000
$ws.Cells.Item(3, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(3, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(3, 16).Value = 1667
$ws.Cells.Item(3, 17).Value = 15
$ws.Cells.Item(3, 18).Value = "Hortaliza"
